$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
for ($r = 13; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    Write-Host $r ":" $cell.Value()
}
Write-Host "UsedRange rows:" $ws.UsedRange.Rows.Count()
Write-Host "UsedRange cols:" $ws.UsedRange.Columns.Count()
